$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the data set.
# Row 7 "Methow River Fawn 02" and row 3 "Entiat River Mills 05".
# Delete bottom-up so row indices of earlier rows stay valid.
$ws.Rows.Item(7).EntireRow.Delete()
$ws.Rows.Item(3).EntireRow.Delete()

# After the two deletions, the former rows 14-17 (Salmon 16-11, Salmon 16-6,
# Salmon 16-9, Tonasket 16-2) are now rows 12-15. Their
# "Riparian-Disturbance_score" (Q) values move to "Riparian-CanopyCover_score"
# (P), and a new "FloodplainConnectivity_score" (M) value of 5 is recorded.
for ($r = 12; $r -le 15; $r++) {
    $q = $ws.Cells.Item($r, 17).Value2
    $ws.Cells.Item($r, 16).Value = $q
    $ws.Cells.Item($r, 17).Value = $null
    $ws.Cells.Item($r, 13).Value = 5
}

# Tonasket 16-2 (now row 15) also had its Flow-SummerBaseFlow_score (L)
# corrected from 0 to 3.
$ws.Cells.Item(15, 12).Value = 3
